# Remove the stale "M2Doc version mismatch" warning that was injected
# between "Basic " and "let" in the first paragraph, restoring the
# original "Basic let demonstration :" text.
#
# The warning is made up of 4 runs sitting right after the "Basic " run
# and right before the "let" run:
#   1) "    "                                   (4 spaces)
#   2) "<---"                                    (orange, highlighted)
#   3) "M2Doc version mismatch: ..."              (orange, highlighted)
#   4) "    "                                   (4 spaces)
#
# Both the "Basic " run and the "let" run carry no explicit run
# formatting, so a plain Range.Delete() spanning the warning would make
# Word coalesce "Basic " and "let" into a single run once they become
# adjacent. To keep them as two distinct runs (matching the original
# document structure / rsid on the "let" run) we briefly tag the "let"
# run with a throw-away character-formatting flag before deleting the
# warning, then clear that flag again -- this stops the engine from
# merging the two runs together.

$d = $word.ActiveDocument

# Locate "Basic " and "let" precisely via Find so we do not depend on
# hard-coded character offsets.
$basicRange = $d.Content
$basicRange.Find.Execute("Basic ") | Out-Null
$afterBasic = $basicRange.End

$letRange = $d.Content
$letRange.Find.Execute("let") | Out-Null
$letStart = $letRange.Start
$letEnd = $letRange.End

# Temporarily mark the "let" run so it keeps its own identity once the
# text in between is removed.
$letMarker = $d.Range($letStart, $letEnd)
$letMarker.Font.Bold = $true

# Delete everything between "Basic " and "let" (the stale warning runs).
$warning = $d.Range($afterBasic, $letStart)
$warning.Delete()

# "let" now immediately follows "Basic " (positions shifted back by the
# length of the deleted warning); recompute its range and undo the
# temporary marker.
$newLetStart = $afterBasic
$newLetEnd = $afterBasic + ($letEnd - $letStart)
$d.Range($newLetStart, $newLetEnd).Font.Bold = $false
